$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Capitan" (AQ) column values for all player rows (2-21).
$ws.Range("AQ2:AQ21").ClearContents()

# Fix accentuation: "Si" -> "Sí" for the substitute players (rows 13-21).
$ws.Range("F13").Value = "Sí"
$ws.Range("F14").Value = "Sí"
$ws.Range("F15").Value = "Sí"
$ws.Range("F16").Value = "Sí"
$ws.Range("F17").Value = "Sí"
$ws.Range("F18").Value = "Sí"
$ws.Range("F19").Value = "Sí"
$ws.Range("F20").Value = "Sí"
$ws.Range("F21").Value = "Sí"
